$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 189 (pushes existing rows 189..259 down to 190..260),
# picking up formatting (incl. the date style on column D) from the row above.
$ws.Rows.Item(189).Insert()

# Columns that stay constant for every "Feria Lagunitas de Puerto Montt" /
# "Perejil" record in this block - copy them straight from the row below
# (the old row 189, now shifted to row 190).
$ws.Range("A189").Value = 4
$ws.Range("B189").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C189").Value = "Los Lagos"
$ws.Range("E189").Value = 10
$ws.Range("F189").Value = 100112044
$ws.Range("G189").Value = "Perejil"
$ws.Range("H189").Value = "Sin especificar"
$ws.Range("I189").Value = "Primera"
$ws.Range("R189").Value = "Hortaliza"

# New record's own data.
$ws.Range("D189").Value = 44726
$ws.Range("J189").Value = 160
$ws.Range("K189").Value = 5000
$ws.Range("L189").Value = 5000
$ws.Range("M189").Value = 5000
$ws.Range("N189").Value = '$/docena de atados (3 kilos)'
$ws.Range("O189").Value = "Región Metropolitana"
$ws.Range("P189").Value = 1667
$ws.Range("Q189").Value = 3
